$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: the old "NIS" id column becomes "username" (A1), the rest
# of the header row stays the same.
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "role"
$ws.Range("D1").Value = "is_active"

# Data row 2: NIS -> new id, role "siswa" -> "wali"
$ws.Range("A2").Value = 888
$ws.Range("B2").Value = "walisantri123"
$ws.Range("C2").Value = "wali"
$ws.Range("D2").Value = $true

# Remove the old second data row entirely (table shrinks from 3 to 2 rows)
$ws.Range("A3:D3").Delete()

# Column widths to match the new content (closest achievable to the
# original author's autofit widths of 11.42578125 / 10.42578125)
$ws.Columns("A").ColumnWidth = 10.65
$ws.Columns("C").ColumnWidth = 9.65

# Selection moves to F4
$ws.Range("F4").Select()
